# Commit: "warning messages ignored in ploting"
#
# The EMISSION sheet's CO2 emission-factor row labels used an inconsistent
# naming scheme (e.g. "1_x_coal_thermal_CO2", "1_5_lignite_CO2", ...) that
# didn't line up with the zero-padded fuel codes used elsewhere in the
# workbook (e.g. "POW_01_x_thermal_coal", "POW_01_05_lignite" on the FUEL
# sheet). That mismatch was causing warning messages to be raised by the
# downstream plotting code, since those labels are used to look up CO2
# factors per fuel. Re-label them to the zero-padded convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EMISSION")

$ws.Range("A2").Value = "01_x_thermal_coal_CO2"
$ws.Range("A3").Value = "01_05_lignite_CO2"
$ws.Range("A4").Value = "02_coal_products_CO2"
$ws.Range("A5").Value = "07_07_gas_diesel_oil_CO2"
$ws.Range("A6").Value = "07_08_fuel_oil_CO2"
$ws.Range("A7").Value = "08_01_natural_gas_CO2"
$ws.Range("A8").Value = "08_01_natural_gas_CCS_CO2"

# The author ended their session with the EMISSION sheet active (instead of
# POWERPLANT), with cell B17 selected.
$ws.Activate()
$ws.Range("B17").Select()
